$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2-81).
# Update it from 45204 (2023-10-05) to 45207 (2023-10-08) for all rows.
$ws.Range("C2:C81").Value = 45207
